$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 2).Value = 238
}
